$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the note text in H10 (shared string for row 10)
$ws.Range("H10").Value = "Performance effect scores of 0 are ranked 5, oroginally ony had upper at 0.1, but multiple were above 0.1"

# Add header-like notes in column H for rows 2-9
# Note: "for Performance Effect scores" must be added to the shared string
# table before "for factorWeight scores" to match the expected string order,
# so set H8 first.
$ws.Range("H8").Value = "for Performance Effect scores"
$ws.Range("H9").Value = "for Performance Effect scores"
$ws.Range("H2").Value = "for factorWeight scores"
$ws.Range("H3").Value = "for factorWeight scores"
$ws.Range("H4").Value = "for factorWeight scores"
$ws.Range("H5").Value = "for factorWeight scores"
$ws.Range("H6").Value = "for factorWeight scores"
$ws.Range("H7").Value = "for factorWeight scores"

# Update F10 value from 0.1 to 10
$ws.Range("F10").Value = 10

# Update selection to H8
$ws.Range("H8").Select()

$wb.Save()
